$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 361, shifting all rows from 361 down by one.
$ws.Rows.Item(361).Insert()

# Populate the newly inserted row 361 with the new weekly data point.
$ws.Cells.Item(361, 1).Value = 6
$ws.Cells.Item(361, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(361, 3).Value = "Metropolitana"
$ws.Cells.Item(361, 4).Value = 44551
$ws.Cells.Item(361, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(361, 5).Value = 13
$ws.Cells.Item(361, 6).Value = 100112012
$ws.Cells.Item(361, 7).Value = "Espinaca"
$ws.Cells.Item(361, 8).Value = "Sin especificar"
$ws.Cells.Item(361, 9).Value = "Primera"
$ws.Cells.Item(361, 10).Value = 580
$ws.Cells.Item(361, 11).Value = 6000
$ws.Cells.Item(361, 12).Value = 6500
$ws.Cells.Item(361, 13).Value = 6216
$ws.Cells.Item(361, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(361, 15).Value = "Región Metropolitana"
$ws.Cells.Item(361, 16).Value = 622
$ws.Cells.Item(361, 17).Value = 10
$ws.Cells.Item(361, 18).Value = "Hortaliza"
